{"js": "// Replace each three-digit-by-one-digit multiplication prompt in the\n// practice table with its new value. Every \"old\" string below is unique\n// within the document, so an exact, case-sensitive whole-string search\n// safely targets exactly one run's text per pair.\nconst replacements = [\n  [\"365\u00d78=\", \"191\u00d75=\"],\n  [\"752\u00d74=\", \"976\u00d72=\"],\n  [\"846\u00d77=\", \"137\u00d75=\"],\n  [\"141\u00d79=\", \"872\u00d75=\"],\n  [\"608\u00d77=\", \"296\u00d74=\"],\n  [\"852\u00d78=\", \"612\u00d73=\"],\n  [\"569\u00d76=\", \"857\u00d74=\"],\n  [\"625\u00d74=\", \"455\u00d79=\"],\n  [\"151\u00d77=\", \"820\u00d73=\"],\n  [\"123\u00d75=\", \"670\u00d75=\"],\n  [\"705\u00d74=\", \"559\u00d73=\"],\n  [\"786\u00d76=\", \"473\u00d76=\"],\n  [\"441\u00d78=\", \"192\u00d74=\"],\n  [\"923\u00d74=\", \"314\u00d76=\"],\n  [\"506\u00d74=\", \"319\u00d77=\"],\n  [\"448\u00d75=\", \"508\u00d77=\"],\n  [\"504\u00d78=\", \"888\u00d79=\"],\n  [\"426\u00d72=\", \"961\u00d78=\"],\n  [\"474\u00d73=\", \"554\u00d77=\"],\n  [\"828\u00d77=\", \"765\u00d72=\"],\n  [\"536\u00d78=\", \"783\u00d76=\"],\n  [\"911\u00d72=\", \"530\u00d76=\"],\n  [\"929\u00d77=\", \"117\u00d73=\"],\n  [\"133\u00d75=\", \"819\u00d78=\"],\n  [\"546\u00d78=\", \"738\u00d79=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, {\n    matchCase: true,\n    matchWholeWord: false,\n  });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each three-digit-by-one-digit multiplication prompt in the\n# practice table with its new value. Every \"old\" string is unique within\n# the document, so Find/Replace with MatchWholeWord off and MatchCase on\n# safely targets exactly one run's text per pair, without touching the\n# trailing \"=\" or any other content.\n$d = $word.ActiveDocument\n\n$replacements = @(\n  @(\"365\u00d78=\", \"191\u00d75=\"),\n  @(\"752\u00d74=\", \"976\u00d72=\"),\n  @(\"846\u00d77=\", \"137\u00d75=\"),\n  @(\"141\u00d79=\", \"872\u00d75=\"),\n  @(\"608\u00d77=\", \"296\u00d74=\"),\n  @(\"852\u00d78=\", \"612\u00d73=\"),\n  @(\"569\u00d76=\", \"857\u00d74=\"),\n  @(\"625\u00d74=\", \"455\u00d79=\"),\n  @(\"151\u00d77=\", \"820\u00d73=\"),\n  @(\"123\u00d75=\", \"670\u00d75=\"),\n  @(\"705\u00d74=\", \"559\u00d73=\"),\n  @(\"786\u00d76=\", \"473\u00d76=\"),\n  @(\"441\u00d78=\", \"192\u00d74=\"),\n  @(\"923\u00d74=\", \"314\u00d76=\"),\n  @(\"506\u00d74=\", \"319\u00d77=\"),\n  @(\"448\u00d75=\", \"508\u00d77=\"),\n  @(\"504\u00d78=\", \"888\u00d79=\"),\n  @(\"426\u00d72=\", \"961\u00d78=\"),\n  @(\"474\u00d73=\", \"554\u00d77=\"),\n  @(\"828\u00d77=\", \"765\u00d72=\"),\n  @(\"536\u00d78=\", \"783\u00d76=\"),\n  @(\"911\u00d72=\", \"530\u00d76=\"),\n  @(\"929\u00d77=\", \"117\u00d73=\"),\n  @(\"133\u00d75=\", \"819\u00d78=\"),\n  @(\"546\u00d78=\", \"738\u00d79=\"),\n)\n\nforeach ($pair in $replacements) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n\n  $range = $d.Content\n  $find = $range.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $oldText\n  $find.Replacement.Text = $newText\n  $find.Forward = $true\n  $find.Wrap = 1\n  $find.Format = $false\n  $find.MatchCase = $true\n  $find.MatchWholeWord = $false\n  $find.MatchWildcards = $false\n  $find.MatchSoundsLike = $false\n  $find.MatchAllWordForms = $false\n\n  $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n"}
